$p = $ppt.ActivePresentation

# Delete the second slide (sldId 260 / rId4); the deck is left with a
# single slide.
$p.Slides.Item(2).Delete()

# On the remaining slide, change the text of the "Step 5" flow-chart box
# from "Step 5: Test Fixed Program" to "Step 5: Test the Fixed Program",
# without disturbing the run/paragraph formatting around it. The shape
# lives one level deep inside a group shape, so both the top level and
# one level of group nesting are checked (this deck never nests deeper
# than that). Recursive helper functions are avoided here because this
# runtime's interpreter does not reliably keep loop-variable state
# isolated across recursive calls.

$s = $p.Slides.Item(1)

$oldText = "Step 5: Test Fixed Program"
$newText = "Step 5: Test the Fixed Program"
$oldLen = $oldText.Length

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shape = $s.Shapes.Item($i)

    if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
        $tr = $shape.TextFrame.TextRange
        $full = $tr.Text
        if ($full.Contains($oldText)) {
            $idx = $full.IndexOf($oldText)
            $start = $idx + 1
            $sub = $tr.Characters($start, $oldLen)
            $sub.Text = $newText
        }
    }

    if ($shape.Type -eq 6) {
        $items = $shape.GroupItems
        for ($j = 1; $j -le $items.Count; $j++) {
            $inner = $items.Item($j)
            if ($inner.HasTextFrame -and $inner.TextFrame.HasText) {
                $itr = $inner.TextFrame.TextRange
                $ifull = $itr.Text
                if ($ifull.Contains($oldText)) {
                    $iidx = $ifull.IndexOf($oldText)
                    $istart = $iidx + 1
                    $isub = $itr.Characters($istart, $oldLen)
                    $isub.Text = $newText
                }
            }
        }
    }
}
